# Weekly update: a new price record (row) is inserted at the top of the
# data block (row 28), pushing all the existing records down by one row.
# The last existing record (old row 79) ends up as the new row 80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 28; everything from 28..79 shifts to 29..80.
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new weekly record.
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44848
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = "Arveja Verde"
$ws.Range("H28").Value = "Perfection"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 26000
$ws.Range("M28").Value = 25500
$ws.Range("N28").Value = "$/malla 25 kilos"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 1020
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
